$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Data rows 16-23: LOXY (CC) / JENNIFFER (PPT, new worker) rows interleaved per period,
# periods now run 2505-2508 instead of 2504-2507 (row 16's period label changes from
# 2507 to 2505, even though B16:D16/F16:G16 keep their original values).

$ws.Range("E16").Value = "2505"

$ws.Range("B17").Value = "PPT"
$ws.Range("C17").Value = "5064239"
$ws.Range("D17").Value = "JENNIFFER PAOLA GARCIA FREITEZ"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 56920
$ws.Range("G17").Value = 1423000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1007027292"
$ws.Range("D18").Value = "LOXY ENRIQUE AVILA ORTIZ"
$ws.Range("E18").Value = "2506"
$ws.Range("F18").Value = 66000
$ws.Range("G18").Value = 1650000

$ws.Range("B19").Value = "PPT"
$ws.Range("C19").Value = "5064239"
$ws.Range("D19").Value = "JENNIFFER PAOLA GARCIA FREITEZ"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 56920
$ws.Range("G19").Value = 1423000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1007027292"
$ws.Range("D20").Value = "LOXY ENRIQUE AVILA ORTIZ"
$ws.Range("E20").Value = "2507"
$ws.Range("F20").Value = 66000
$ws.Range("G20").Value = 1650000

$ws.Range("B21").Value = "PPT"
$ws.Range("C21").Value = "5064239"
$ws.Range("D21").Value = "JENNIFFER PAOLA GARCIA FREITEZ"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 56920
$ws.Range("G21").Value = 1423000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1007027292"
$ws.Range("D22").Value = "LOXY ENRIQUE AVILA ORTIZ"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 66000
$ws.Range("G22").Value = 1650000

$ws.Range("B23").Value = "PPT"
$ws.Range("C23").Value = "5064239"
$ws.Range("D23").Value = "JENNIFFER PAOLA GARCIA FREITEZ"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56920
$ws.Range("G23").Value = 1423000
